# Generate Report for Handback
# Adds a new handback entry (file 2b625ea6-3f28-4cf2-b7b3-881c24aa0c3c.md) as
# row 3 on each of the three report sheets (Overview, zh-cn, de-de), mirroring
# the existing row 2 entry for e414559a-85d2-4c60-8b29-5c9aa639a168.md /
# 23cb1c5e-30bc-4c0f-add8-5146f899bdb5.md.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Update the "current" handback file's identity + timestamp everywhere it's
# referenced (rename e414559a-... -> 23cb1c5e-... and bump the handoff time).
# ---------------------------------------------------------------------------
$oldGuid = "e414559a-85d2-4c60-8b29-5c9aa639a168"
$newGuid = "23cb1c5e-30bc-4c0f-add8-5146f899bdb5"
$secondGuid = "2b625ea6-3f28-4cf2-b7b3-881c24aa0c3c"

$oldZhHash = "3bfab2c72810c51a52d881e839aa9c9ddc79520a"
$newZhHash = "13bca1b8e7a2e7d608916f9535620632f3c2f157"
$secondZhHash = "feb03cd2ba7c5581258dc19bb7cc1ff745ff7a33"

$oldDeHash = "3bfab2c72810c51a52d881e839aa9c9ddc79520a"
$newDeHash = "13bca1b8e7a2e7d608916f9535620632f3c2f157"
$secondDeHash = "feb03cd2ba7c5581258dc19bb7cc1ff745ff7a33"

# ===========================================================================
# Sheet "Overview"
# ===========================================================================
$wsOv = $wb.Worksheets.Item("Overview")

# Rename the existing (row 2) file references + refresh its generate date.
$wsOv.Range("A2").Value = "$newGuid.md"
$wsOv.Range("B2").Value = "e2e\$newGuid.md"
$wsOv.Range("G2").Value = "2016-09-06 21:22:52"
$wsOv.Range("G2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOv.Hyperlinks.Item(1).Address = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d76fb1add516af8839da5f2145d6e4ecc92013b6/e2e/$newGuid.md"
$wsOv.Hyperlinks.Item(1).TextToDisplay = "e2e\$newGuid.md"

# Add the new (row 3) entry for the second handback file.
$loOv = $wsOv.ListObjects.Item("Overview")
$loOv.ListRows.Add() | Out-Null

$wsOv.Range("A3").Value = "$secondGuid.md"
$wsOv.Range("B3").Value = "e2e\$secondGuid.md"
$wsOv.Range("C3").Value = ".md"
$wsOv.Range("E3").Value = "Handed back: in sync with en-US"
$wsOv.Range("F3").Value = "Handed back: in sync with en-US"
$wsOv.Range("G3").Value = "2016-09-06 21:22:52"
$wsOv.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOv.Hyperlinks.Add($wsOv.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d76fb1add516af8839da5f2145d6e4ecc92013b6/e2e/$secondGuid.md", [Type]::Missing, [Type]::Missing, "e2e\$secondGuid.md") | Out-Null

# ===========================================================================
# Sheet "zh-cn"
# ===========================================================================
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Range("G2").Value = "$newGuid.$newZhHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-09-06 21:22:46"
$wsZh.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I2").Value = "$newGuid.md"
$wsZh.Range("J2").Value = "$newGuid.$newZhHash.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-09-06 21:23:23"
$wsZh.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZh.Hyperlinks.Item(1).Address = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d76fb1add516af8839da5f2145d6e4ecc92013b6/e2e/$newGuid.md"
$wsZh.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"
$wsZh.Hyperlinks.Item(2).Address = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/71c6288f3a11a64220f73ce7f8f11a4b6602b59a/e2e/$newGuid.md"
$wsZh.Hyperlinks.Item(2).TextToDisplay = "$newGuid.md"

$loZh = $wsZh.ListObjects.Item("zh-cn")
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A3").Value = "$secondGuid.md"
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "True"
$wsZh.Range("G3").Value = "$secondGuid.$secondZhHash.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-09-06 21:22:46"
$wsZh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I3").Value = "$secondGuid.md"
$wsZh.Range("J3").Value = "$secondGuid.$secondZhHash.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-09-06 21:23:23"
$wsZh.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L3").Value = ""
$wsZh.Range("M3").Value = "True"
$wsZh.Range("N3").Value = ""
$wsZh.Range("O3").Value = "False"
$wsZh.Range("P3").Value = ""

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d76fb1add516af8839da5f2145d6e4ecc92013b6/e2e/$secondGuid.md", [Type]::Missing, [Type]::Missing, "$secondGuid.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/71c6288f3a11a64220f73ce7f8f11a4b6602b59a/e2e/$secondGuid.md", [Type]::Missing, [Type]::Missing, "$secondGuid.md") | Out-Null

# ===========================================================================
# Sheet "de-de"
# ===========================================================================
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Range("G2").Value = "$newGuid.$newDeHash.de-de.xlf"
$wsDe.Range("J2").Value = "$newGuid.$newDeHash.de-de.xlf"
$wsDe.Range("K2").Value = "2016-09-06 21:23:31"
$wsDe.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I2").Value = "$newGuid.md"

$wsDe.Hyperlinks.Item(1).Address = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d76fb1add516af8839da5f2145d6e4ecc92013b6/e2e/$newGuid.md"
$wsDe.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"
$wsDe.Hyperlinks.Item(2).Address = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/66e30d5bed1cae642846f46f4e4a65ab872b9270/e2e/$newGuid.md"
$wsDe.Hyperlinks.Item(2).TextToDisplay = "$newGuid.md"

$loDe = $wsDe.ListObjects.Item("de-de")
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A3").Value = "$secondGuid.md"
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "True"
$wsDe.Range("G3").Value = "$secondGuid.$secondDeHash.de-de.xlf"
$wsDe.Range("H3").Value = "2016-09-06 21:22:52"
$wsDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I3").Value = "$secondGuid.md"
$wsDe.Range("J3").Value = "$secondGuid.$secondDeHash.de-de.xlf"
$wsDe.Range("K3").Value = "2016-09-06 21:23:31"
$wsDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L3").Value = ""
$wsDe.Range("M3").Value = "True"
$wsDe.Range("N3").Value = ""
$wsDe.Range("O3").Value = "False"
$wsDe.Range("P3").Value = ""

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d76fb1add516af8839da5f2145d6e4ecc92013b6/e2e/$secondGuid.md", [Type]::Missing, [Type]::Missing, "$secondGuid.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/66e30d5bed1cae642846f46f4e4a65ab872b9270/e2e/$secondGuid.md", [Type]::Missing, [Type]::Missing, "$secondGuid.md") | Out-Null
